$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the old "rbf(size=46)" column (old E) into the new F column.
$ws.Range("F8").Value = $ws.Range("E8").Value()
$ws.Range("F9").Value = $ws.Range("E9").Value()
$ws.Range("F10").Value = $ws.Range("E10").Value()
$ws.Range("F11").Value = $ws.Range("E11").Value()

$ws.Range("F8").Style = $ws.Range("D8").Style()

# Clear the old E11 text (moved to F11) - will be overwritten below (blank)
$ws.Range("E11").ClearContents()

# New rbf(size=8, rbfweights) description in E8, with its own font/wrap style.
$ws.Range("E8").Value = "rbf(size=8, rbfweights)
rbf(nn_df, nn_df`$shares, size=8, linOut=TRUE,
                 initFunc = ""RBF_Weights"", initFuncParams = c(0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5, 0.5),
                 learnFunc = ""RadialBasisLearning"", learnFuncParams = c(1e-05, 1e-05, 1e-05, 1e-05, 1e-05, 1e-05, 1e-05,1e-05))"

$ws.Range("E9").Value = 0.02419791
$ws.Range("E10").Value = 0.026383509999999999

# D10 is no longer bold-red; make it the plain bold style (same as A9/A10/B2 header cells).
$ws.Range("D10").Style = $ws.Range("A10").Style()

Write-Output "done"
